# Update SA templates with new field:
# Insert a new "Ethnic or Racial Identity" field group (3 columns) right
# after the existing "Gender Identity, Sexuality" field group (before the
# "Pronouns" group), matching the repeating
# (open-description-tag, bold label, close-description-tag) column pattern
# used throughout row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank columns at J:L, shifting everything from the old
# column J onward to the right by 3 columns.
$ws.Range("J1:L1").EntireColumn.Insert() | Out-Null

# Populate the new field group.
$ws.Range("J1").Value = "<mods:description>Ethnic or Racial Identity: "
$ws.Range("K1").Value = "Ethnic or Racial Identity"
$ws.Range("K1").Font.Bold = $true
$ws.Range("L1").Value = "</mods:description>"

# Match the author's final selection/cursor position.
$ws.Range("J1").Select() | Out-Null
